# Strip the trailing footnote markers (e.g. " [1]", " [5]") from vaccine/
# brand-name labels across all four sheets, and collapse the handful of
# labels that were wrapped onto a second line (embedded newline) into a
# single line with a space instead. Also folds the lone "Afluria\nQuadrivalent"
# label (Adult Influenza sheet, B9/B10) into the already-existing single-line
# "Afluria Quadrivalent" text used elsewhere on that same sheet.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# --- Sheet 1 ---
$ws1.Range("A2").Value = 'DTaP '
$ws1.Range("A3").Value = 'DTaP '
$ws1.Range("A4").Value = 'DTaP-IPV '
$ws1.Range("A5").Value = 'DTaP-IPV '
$ws1.Range("A6").Value = 'DTaP-IPV '
$ws1.Range("A7").Value = 'DTaP-Hep B-IPV '
$ws1.Range("A8").Value = 'DTaP-IP-HI '
$ws1.Range("A9").Value = 'e-IPV '
$ws1.Range("A10").Value = 'Hepatitis A Pediatric '
$ws1.Range("A11").Value = 'Hepatitis A Pediatric '
$ws1.Range("A12").Value = 'Hepatitis A-Hepatitis B 18 only '
$ws1.Range("A13").Value = 'Hepatitis B  Pediatric/Adolescent'
$ws1.Range("A14").Value = 'Hepatitis B  Pediatric/Adolescent'
$ws1.Range("B14").Value = 'Recombivax HB'
$ws1.Range("A15").Value = 'Hib '
$ws1.Range("A16").Value = 'Hib '
$ws1.Range("A17").Value = 'Hib '
$ws1.Range("A18").Value = 'HPV - Human Papillomavirus 9-valent '
$ws1.Range("A19").Value = 'MENB - Meningococcal Group B '
$ws1.Range("A20").Value = 'MENB - Meningococcal Group B '
$ws1.Range("A21").Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws1.Range("A22").Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws1.Range("A23").Value = 'Measles, Mumps and Rubella (MMR) '
$ws1.Range("A24").Value = 'MMR/Varicella '
$ws1.Range("A25").Value = 'Pneumococcal 13-valent  (Pediatric)'
$ws1.Range("A27").Value = 'Rotavirus, Live, Oral, Pentavalent '
$ws1.Range("A28").Value = 'Rotavirus, Live, Oral, Pentavalent '
$ws1.Range("A29").Value = 'Rotavirus, Live, Oral, Oral '
$ws1.Range("A30").Value = 'Tetanus and Diphtheria Toxoids '
$ws1.Range("A31").Value = 'Tetanus and Diphtheria Toxoids '
$ws1.Range("A32").Value = 'Tetanus and Diphtheria Toxoids '
$ws1.Range("A33").Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws1.Range("A34").Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws1.Range("A35").Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws1.Range("A36").Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws1.Range("A37").Value = 'Varicella '

# --- Sheet 2 ---
$ws2.Range("A2").Value = 'Hepatitis A Adult '
$ws2.Range("A3").Value = 'Hepatitis A Adult '
$ws2.Range("A4").Value = 'Hepatitis A Adult '
$ws2.Range("A5").Value = 'Hepatitis A-Hepatitis B Adult '
$ws2.Range("A6").Value = 'Hepatitis B Adult '
$ws2.Range("A7").Value = 'Hepatitis B Adult '
$ws2.Range("A8").Value = 'Hepatitis B Adult '
$ws2.Range("A9").Value = 'HPV-Human Papillomavirus 9 Valent '
$ws2.Range("A10").Value = 'Measles, Mumps,  Rubella '
$ws2.Range("A11").Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws2.Range("A12").Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws2.Range("A13").Value = 'MENB - Meningococcal Group B '
$ws2.Range("A14").Value = 'MENB - Meningococcal Group B '
$ws2.Range("A15").Value = 'Pneumococcal 13-valent '
$ws2.Range("A17").Value = 'Tetanus and Diphtheria Toxoids '
$ws2.Range("A18").Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws2.Range("A19").Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws2.Range("A20").Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws2.Range("A21").Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws2.Range("A22").Value = 'Varicella '

# --- Sheet 3 ---
$ws3.Range("A2").Value = 'Influenza  (Age 6 months and older)'
$ws3.Range("B2").Value = 'Fluzone Quadrivalent'
$ws3.Range("A3").Value = 'Influenza  (Age 6 months and older)'
$ws3.Range("B3").Value = 'Fluzone Quadrivalent'
$ws3.Range("A4").Value = 'Influenza  (Age 6 months and older)'
$ws3.Range("B4").Value = 'Fluzone Quadrivalent'
$ws3.Range("A5").Value = 'Influenza  (Age 6 months and older)'
$ws3.Range("B5").Value = 'Fluarix Quadrivalent'
$ws3.Range("A6").Value = 'Influenza  (Age 6 months and older)'
$ws3.Range("B6").Value = 'FluLaval Quadrivalent'
$ws3.Range("A7").Value = 'Influenza  (Age 4 years and older)'
$ws3.Range("A8").Value = 'Influenza  (Age 4 years and older)'
$ws3.Range("A9").Value = 'Influenza  (Age 6 -35 months)'
$ws3.Range("A10").Value = 'Influenza  (Age 36 months and older)'
$ws3.Range("A11").Value = 'Influenza  (Age 6 months and older)'
$ws3.Range("A12").Value = 'Influenza  Live, Intranasal (Age 2-49 years)'
$ws3.Range("B12").Value = 'FluMist Quadrivalent'

# --- Sheet 4 ---
$ws4.Range("A2").Value = 'Influenza  (Age 6 months and older)'
$ws4.Range("B2").Value = 'Fluzone Quadrivalent'
$ws4.Range("A3").Value = 'Influenza  (Age 6 months and older)'
$ws4.Range("B3").Value = 'Fluzone Quadrivalent'
$ws4.Range("A4").Value = 'Influenza  (Age 6 months and older)'
$ws4.Range("B4").Value = 'Fluzone Quadrivalent'
$ws4.Range("A5").Value = 'Influenza  (Age 6 months and older)'
$ws4.Range("B5").Value = 'Fluarix Quadrivalent'
$ws4.Range("A6").Value = 'Influenza  (Age 6 months and older)'
$ws4.Range("B6").Value = 'FluLaval Quadrivalent'
$ws4.Range("A7").Value = 'Influenza  (Age 4 years and older)'
$ws4.Range("A8").Value = 'Influenza  (Age 4 years and older)'
$ws4.Range("A9").Value = 'Influenza  (Age 36 months and older)'
$ws4.Range("B9").Value = 'Afluria Quadrivalent'
$ws4.Range("A10").Value = 'Influenza  (Age 6 months and older)'
$ws4.Range("B10").Value = 'Afluria Quadrivalent'
$ws4.Range("A11").Value = 'Influenza  Live, Intranasal (Age 2-49 years)'
$ws4.Range("B11").Value = 'FluMist Quadrivalent'
